$wb = $excel.ActiveWorkbook

# Rename sheet "TestControl" -> "LoginData"
$ws = $wb.Worksheets.Item("TestControl")
$ws.Name = "LoginData"

$ws.Range("A1").Value = "StartLoginTest"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Pasword"
$ws.Range("D1").Value = "RanMode"

$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "stasman9@mail.ru"
$ws.Range("C2").Value = "123456789@"
$ws.Range("D2").Value = "Y"

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "stasman9@mail.ru"
$ws.Range("C3").Value = "123456789@"
$ws.Range("D3").Value = "Y"

$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "stasman9@mail.ru"
$ws.Range("C4").Value = "123456789@"
$ws.Range("D4").Value = "Y"

Write-Host "done"
